$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "1차 선형 미분 방정식의 해법"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2021/05/08/first_order_linear_equations.html"

$ws.Range("D8").Value = "제주어 기계번역 모델과 음성합성 모델에 관한 연구를 소개합니다."

$ws.Range("D20").Value = "[파이썬 간단한 게임 만들기] 9. 오목 아니고 4목"
$ws.Range("E20").Value = "https://ai-creator.tistory.com/537"

$ws.Range("D51").Value = "[python] join 함수로 리스트의 요소들 하나의 문자열로 합치기"
$ws.Range("E51").Value = "https://bskyvision.com/1189"
